$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# 1. Change the SimulId used for the first train/infer row
$ws.Range("A6").Value = 9044

# 2. Change the Engine used for the second train row
$ws.Range("B7").Value = "Engine2"

# 3. Paste the (recalculated) inferred values of E6:H6 as plain values,
#    one per row, into E8:E11
$ws.Range("E6").Copy()
$ws.Range("E8").PasteSpecial(-4163) | Out-Null
$ws.Range("F6").Copy()
$ws.Range("E9").PasteSpecial(-4163) | Out-Null
$ws.Range("G6").Copy()
$ws.Range("E10").PasteSpecial(-4163) | Out-Null
$ws.Range("H6").Copy()
$ws.Range("E11").PasteSpecial(-4163) | Out-Null

$ws.Range("D7").Select()
